$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Bool Column" header in E1
$ws.Range("E1").Value = "Bool Column"

# Row 2: real boolean TRUE
$ws.Range("E2").Value = $true

# Row 3: real boolean FALSE
$ws.Range("E3").Value = $false

# Row 4: plain numeric 1 (schema override -> stored as number, not boolean)
$ws.Range("E4").Value = 1

# Row 5 intentionally left blank (no value) to mirror missing data

# Row 6: plain numeric 0 (schema override -> stored as number, not boolean)
$ws.Range("E6").Value = 0

# Size the new column to fit its content, like the bestFit pass Excel runs
# automatically for a freshly introduced column.
$ws.Columns.Item(5).ColumnWidth = 10.1666666

# Move the active selection to E7, matching the post-edit cursor position
$ws.Range("E7").Select() | Out-Null
